$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns are treated as text so that
# numeric-looking strings (e.g. "1.00", "247.00") are preserved exactly
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '42.726.14'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '2.544.52'
$ws.Range('E3').Value = '  +0.60%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '308.86'
$ws.Range('E5').Value = '  -2.44%  '
$ws.Range('D6').Value = '97.49'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = '0.572'
$ws.Range('E7').Value = '  -0.60%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.529'
$ws.Range('E9').Value = '  -0.77%  '
$ws.Range('D10').Value = '35.53'
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('D11').Value = '0.0806'
$ws.Range('E11').Value = '  -0.58%  '
$ws.Range('D12').Value = '7.40'
$ws.Range('E12').Value = '  -2.84%  '
$ws.Range('D14').Value = '2.934.91'
$ws.Range('E14').Value = '  +0.80%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '15.76'
$ws.Range('E15').Value = '  +4.08%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.570.91'
$ws.Range('E16').Value = '  +2.25%  '
$ws.Range('D17').Value = '0.834'
$ws.Range('E17').Value = '  -1.72%  '
$ws.Range('D18').Value = '42.716.46'
$ws.Range('E18').Value = '  -0.30%  '
$ws.Range('D19').Value = '6.74'
$ws.Range('E19').Value = '  -1.58%  '
$ws.Range('D20').Value = '12.37'
$ws.Range('E20').Value = '  -3.00%  '
$ws.Range('D21').Value = '0.0₃0956'
$ws.Range('E21').Value = '  -0.64%  '
$ws.Range('D22').Value = '69.23'
$ws.Range('E22').Value = '  -0.61%  '
$ws.Range('D23').Value = '247.00'
$ws.Range('E23').Value = '  -2.10%  '
$ws.Range('D24').Value = '2.91'
$ws.Range('E24').Value = '  -1.52%  '
$ws.Range('D25').Value = '2.04'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').Value = '26.54'
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').Value = '2.36'
$ws.Range('E28').Value = '  -2.15%  '
$ws.Range('D29').Value = '40.16'
$ws.Range('E29').Value = '  -2.55%  '
$ws.Range('D30').Value = '10.13'
$ws.Range('E30').Value = '  -2.77%  '
$ws.Range('D31').Value = '157.67'
$ws.Range('E31').Value = '  -1.17%  '
$ws.Range('D32').Value = '5.72'
$ws.Range('E33').Value = '  +0.33%  '
$ws.Range('E34').Value = '  -0.89%  '
$ws.Range('D35').Value = '2.08'
$ws.Range('E35').Value = '  -3.63%  '
$ws.Range('E36').Value = '  -3.62%  '
$ws.Range('D37').Value = '18.44'
$ws.Range('E37').Value = '  -2.39%  '
$ws.Range('D38').Value = '2.57'
$ws.Range('E38').Value = '  +4.76%  '
$ws.Range('E39').Value = '  -1.62%  '
$ws.Range('E40').Value = '  -0.85%  '
$ws.Range('D41').Value = '22.46'
$ws.Range('E41').Value = '  +2.67%  '
$ws.Range('D42').Value = '4.04'
$ws.Range('E42').Value = '  +5.18%  '
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('E44').Value = '  -2.08%  '
$ws.Range('D45').Value = '1.991.82'
$ws.Range('E45').Value = '  -1.49%  '
$ws.Range('D46').Value = '3.20'
$ws.Range('E46').Value = '  -2.46%  '
$ws.Range('D47').Value = '9.03'
$ws.Range('E47').Value = '  -0.66%  '
$ws.Range('D48').Value = '2.789.05'
$ws.Range('E48').Value = '  +0.81%  '
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('D50').Value = '80.66'
$ws.Range('E50').Value = '  -4.13%  '
$ws.Range('D51').Value = '73.39'
$ws.Range('E51').Value = '  -2.61%  '

# Restore the original (default/no explicit number format) style for the
# Price and Volume columns now that the text values have been written.
$ws.Range("D2:E51").Style = "Normal"
